$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (D1:J1), matching style of existing headers (A1:C1)
$headers = @("Migration Date", "FINANCIAL INSTITUTION NAME", "Entity ID", "CURRENT SWITCH", "Old Platform", "New Platform", "Service")
$col = 4
foreach ($h in $headers) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Style = $ws.Cells.Item(1, 1).Style
    $col++
}

# Clear C6:C8 (previously had "NA", now blank)
$ws.Range("C6:C8").Value = ""

# Add blank values in D2:J8 so the cells exist as empty inline strings
$ws.Range("D2:J8").Value = ""

# New row 9 - "client setting/value" mirrored row, with migration block shifted into D:J
$ws.Cells.Item(9, 1).Value = "client setting/value"
$ws.Cells.Item(9, 2).Value = ""
$ws.Cells.Item(9, 3).Value = ""
$ws.Cells.Item(9, 4).Value = "2025-10-16 00:00:00"
$ws.Cells.Item(9, 5).Value = "YYY"
$ws.Cells.Item(9, 6).Value = "123ABX007"
$ws.Cells.Item(9, 7).Value = "FISB"
$ws.Cells.Item(9, 8).Value = "NA"
$ws.Cells.Item(9, 9).Value = "NA"
$ws.Cells.Item(9, 10).Value = "NA"
